$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column (D) for numeric-looking values so Excel
# does not auto-convert them to numbers; we restore the style afterwards so
# no residual style/format difference remains on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '22.423.31'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.567.24'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').Value = '285.18'
$ws.Range('E6').Value = '  -2.30%  '
$ws.Range('D7').Value = '0.3626'
$ws.Range('E7').Value = '  -2.62%  '
$ws.Range('D8').Value = '48.50'
$ws.Range('E8').Value = '  -2.78%  '
$ws.Range('D9').Value = '0.3321'
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('D10').Value = '1.121'
$ws.Range('D11').Value = '0.07390'
$ws.Range('E11').Value = '  -2.29%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('E13').Value = '  -2.32%  '
$ws.Range('D14').Value = '5.941'
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('D15').Value = '6.899'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('D16').Value = '1.567.10'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('E17').Value = '  -1.66%  '
$ws.Range('D18').Value = '87.96'
$ws.Range('E18').Value = '  -3.32%  '
$ws.Range('D19').Value = '0.06698'
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = '6.326'
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('D24').Value = '22.423.73'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').Value = '2.374'
$ws.Range('E25').Value = '  +1.56%  '
$ws.Range('D26').Value = '2.537'
$ws.Range('E26').Value = '  -5.26%  '
$ws.Range('D27').Value = '150.39'
$ws.Range('E27').Value = '  +1.28%  '
$ws.Range('D28').Value = '19.37'
$ws.Range('E28').Value = '  -3.56%  '
$ws.Range('D29').Value = '4.993'
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('D30').Value = '123.70'
$ws.Range('D31').Value = '1.743.63'
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('D32').Value = '1.036'
$ws.Range('E32').Value = '  -2.07%  '
$ws.Range('D33').Value = '2.005'
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('D34').Value = '6.087'
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('D35').Value = '9.809'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('D36').Value = '0.08233'
$ws.Range('E36').Value = '  -1.84%  '
$ws.Range('D37').Value = '0.02409'
$ws.Range('E37').Value = '  -2.68%  '
$ws.Range('D38').Value = '0.2233'
$ws.Range('E38').Value = '  -3.10%  '
$ws.Range('D39').Value = '0.06421'
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('D40').Value = '5.363'
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('D41').Value = '1.287'
$ws.Range('E41').Value = '  -5.09%  '
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('D43').Value = '11.17'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '13.73'
$ws.Range('E45').Value = '  -1.47%  '
$ws.Range('D46').Value = '0.6042'
$ws.Range('E46').Value = '  +4.06%  '
$ws.Range('D47').Value = '3.751'
$ws.Range('E47').Value = '  -1.69%  '
$ws.Range('D48').Value = '2.027'
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('D49').Value = '123.24'
$ws.Range('E49').Value = '  -5.30%  '
$ws.Range('D50').Value = '1.210'
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('D51').Value = '0.07200'
$ws.Range('E51').Value = '  -1.61%  '

# Restore original (default) style on the price column now that the
# text values have been written, so cell formatting matches the source.
$priceRange.Style = "Normal"

